# Add a new "2021年" data row (row 13) under the existing "引进电子出版物版权总数"
# table on Sheet1, following the same layout as the prior year rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the year label and the numeric figures that are actually reported
# for 2021. Columns with no reported figure get a lone "'" (apostrophe) so
# Excel stores them as empty *text* cells -- matching how the other rows in
# this sheet represent "no data" -- rather than leaving the cell fully
# blank.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = "'"
$ws.Range("C13").Value = 17
$ws.Range("D13").Value = "'"
$ws.Range("E13").Value = "'"
$ws.Range("F13").Value = "'"
$ws.Range("G13").Value = "'"
$ws.Range("H13").Value = 7
$ws.Range("I13").Value = "'"
$ws.Range("J13").Value = "'"
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = "'"
$ws.Range("M13").Value = 2
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 28

# Match the formatting of the row above (bold/centered/bordered label cell
# in column A, plain numeric cells elsewhere) by copying its formats down.
# Doing this *after* setting the values also clears any stray "quote
# prefix" formatting flag left behind by the "'" entries above.
$ws.Range("A12:O12").Copy()
$ws.Range("A13:O13").PasteSpecial(-4122)
